# Updates cryptocurrency price/volume figures in the "cryptos" sheet
# to match the latest scrape (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.039.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.637.84"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "  +0.61%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("E8").Value = "  -1.70%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0628"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.68"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.62%  "
$ws.Range("E11").Value = "  +0.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.714.69"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.83%  "
$ws.Range("E13").Value = "  -1.53%  "
$ws.Range("E14").Value = "  -2.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "62.29"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₃0748"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.056.99"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.41%  "
$ws.Range("E18").Value = "  +0.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "191.31"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.83%  "
$ws.Range("E20").Value = "  -1.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.16"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.58%  "
$ws.Range("E23").Value = "  +1.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "143.76"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.36%  "
$ws.Range("B25").Value = "BinanceUSD"
$ws.Range("C25").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.66%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.77"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.80%  "
$ws.Range("E29").Value = "  -0.23%  "
$ws.Range("E31").Value = "  -2.23%  "
$ws.Range("E32").Value = "  -2.95%  "
$ws.Range("E33").Value = "  -1.09%  "
$ws.Range("E34").Value = "  -0.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.879"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.128.11"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.35%  "
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.528"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.77%  "
$ws.Range("E39").Value = "  -0.85%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "99.08"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.786"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.30"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.11%  "
$ws.Range("E43").Value = "  -0.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "55.60"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.79%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0525"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.65%  "
$ws.Range("E46").Value = "  +1.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.414"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.14%  "
$ws.Range("E48").Value = "  -1.28%  "
$ws.Range("E49").Value = "  +0.35%  "
$ws.Range("E50").Value = "  -2.82%  "
$ws.Range("E51").Value = "  -0.39%  "
